$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded "Gestational Age Units" column (originally column X).
$ws.Columns("X").Delete()

# Remove the trailing, always-empty "Quantity" column (originally column AE,
# now column AD after the deletion above).
$ws.Columns("AD").Delete()

# Rename "Gestational Age" (now column W) to reflect the merged unit info.
$ws.Range("W1").Value = "Gestational Age (Weeks)"
$ws.Range("W1").ColumnWidth = 21.375

# Match the author's final selection in the saved workbook.
$ws.Range("F8").Select()
